{"js": "// Replace the date line and each two-digit multiplication problem in the\n// table with the updated values, per the commit's regenerated content.\nconst replacements = [\n  [\"2024-04-02 Tuesday\", \"2024-04-03 Wednesday\"],\n  [\"88\u00d779=\", \"37\u00d795=\"],\n  [\"91\u00d757=\", \"73\u00d768=\"],\n  [\"56\u00d727=\", \"68\u00d743=\"],\n  [\"33\u00d779=\", \"17\u00d775=\"],\n  [\"34\u00d732=\", \"77\u00d774=\"],\n  [\"52\u00d754=\", \"15\u00d789=\"],\n  [\"74\u00d744=\", \"48\u00d746=\"],\n  [\"93\u00d716=\", \"27\u00d787=\"],\n  [\"63\u00d721=\", \"75\u00d711=\"],\n  [\"26\u00d769=\", \"66\u00d760=\"],\n  [\"18\u00d756=\", \"91\u00d798=\"],\n  [\"23\u00d751=\", \"82\u00d740=\"],\n  [\"93\u00d723=\", \"56\u00d780=\"],\n  [\"69\u00d792=\", \"81\u00d784=\"],\n  [\"19\u00d738=\", \"99\u00d735=\"],\n  [\"48\u00d724=\", \"79\u00d742=\"],\n  [\"97\u00d771=\", \"72\u00d780=\"],\n  [\"33\u00d796=\", \"28\u00d745=\"],\n  [\"37\u00d744=\", \"53\u00d760=\"],\n  [\"32\u00d718=\", \"62\u00d737=\"],\n  [\"22\u00d724=\", \"24\u00d780=\"],\n  [\"51\u00d794=\", \"57\u00d797=\"],\n  [\"86\u00d791=\", \"74\u00d752=\"],\n  [\"45\u00d728=\", \"17\u00d730=\"],\n  [\"12\u00d768=\", \"82\u00d743=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and each two-digit multiplication problem in the\n# table with the regenerated values for this day's worksheet.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-04-02 Tuesday\", \"2024-04-03 Wednesday\"),\n    @(\"88\u00d779=\", \"37\u00d795=\"),\n    @(\"91\u00d757=\", \"73\u00d768=\"),\n    @(\"56\u00d727=\", \"68\u00d743=\"),\n    @(\"33\u00d779=\", \"17\u00d775=\"),\n    @(\"34\u00d732=\", \"77\u00d774=\"),\n    @(\"52\u00d754=\", \"15\u00d789=\"),\n    @(\"74\u00d744=\", \"48\u00d746=\"),\n    @(\"93\u00d716=\", \"27\u00d787=\"),\n    @(\"63\u00d721=\", \"75\u00d711=\"),\n    @(\"26\u00d769=\", \"66\u00d760=\"),\n    @(\"18\u00d756=\", \"91\u00d798=\"),\n    @(\"23\u00d751=\", \"82\u00d740=\"),\n    @(\"93\u00d723=\", \"56\u00d780=\"),\n    @(\"69\u00d792=\", \"81\u00d784=\"),\n    @(\"19\u00d738=\", \"99\u00d735=\"),\n    @(\"48\u00d724=\", \"79\u00d742=\"),\n    @(\"97\u00d771=\", \"72\u00d780=\"),\n    @(\"33\u00d796=\", \"28\u00d745=\"),\n    @(\"37\u00d744=\", \"53\u00d760=\"),\n    @(\"32\u00d718=\", \"62\u00d737=\"),\n    @(\"22\u00d724=\", \"24\u00d780=\"),\n    @(\"51\u00d794=\", \"57\u00d797=\"),\n    @(\"86\u00d791=\", \"74\u00d752=\"),\n    @(\"45\u00d728=\", \"17\u00d730=\"),\n    @(\"12\u00d768=\", \"82\u00d743=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute([ref]$old, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$new, [ref]2)\n}\n"}
